# "Update with new data from Brownie"
# Refresh the raw comparison data on the "copy raw data here" sheet. The
# "Formatted table" sheet pulls every value from here via formulas, so it
# recalculates automatically.

$wb  = $excel.ActiveWorkbook
$raw = $wb.Worksheets.Item("copy raw data here")
$fmt = $wb.Worksheets.Item("Formatted table")

$raw.Range("B2").Value = -29.4
$raw.Range("C2").Value = 0.74349043733634101
$raw.Range("D2").Value = 0.60356346447113396
$raw.Range("E2").Value = 0.928692447243329

$raw.Range("B3").Value = 12.9
$raw.Range("C3").Value = 0.84472743864157596
$raw.Range("D3").Value = 0.62015451485279505
$raw.Range("E3").Value = 1.3660493018355599

$raw.Range("B4").Value = -13.9
$raw.Range("C4").Value = 0.74635829171337797
$raw.Range("D4").Value = 0.6089281564347
$raw.Range("E4").Value = 1.08033256844143

$raw.Range("B5").Value = -40.9
$raw.Range("C5").Value = 0.92046421133088696
$raw.Range("D5").Value = 0.46822645079120701
$raw.Range("E5").Value = 0.80608254529997103

$raw.Range("B6").Value = -24
$raw.Range("C6").Value = 0.78455346348008304
$raw.Range("D6").Value = 0.55925621500540501
$raw.Range("E6").Value = 0.96342100094975203

$raw.Range("B7").Value = 42.7
$raw.Range("C7").Value = 0.87810338680559996
$raw.Range("D7").Value = 0.67054333105544495
$raw.Range("E7").Value = 1.3908322153116499

$raw.Range("B8").Value = -16.6
$raw.Range("C8").Value = 0.73479916717468197
$raw.Range("D8").Value = 0.62411197467359902
$raw.Range("E8").Value = 1.05969561462669

$excel.Calculate()

# Move the lingering cell cursor on the raw-data sheet, then re-activate the
# "Formatted table" tab so it stays the selected/visible sheet, matching the
# saved view state.
$raw.Range("G24").Select() | Out-Null
$fmt.Activate() | Out-Null
